$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Data Source" note to reflect the combined 2017/19 survey wave.
$ws.Range("A3").Value = "Data Source: [2017/19 Puget Sound Regional Household Travel Survey] (https://www.psrc.org/household-travel-survey-program) "

# Reflect the new active selection on the About tab (cell A3).
$ws.Range("A3").Select()
